# "Modified Novice Tester 2 steps"
# Update the NoviceTester2Steps sample row (row 2) so that:
#  - ZipCode (I2) is stored as text "363310" instead of a number
#  - manufacturer (P2) changes from "Samsung" to "S"
#  - model (Q2) changes from "Galaxy" to "GG"
#  - mobile_OS (R2) changes from "Android 4.3" to "A"
#  - the active selection on the sheet moves from R2 to S2
#  - the sheet page setup is switched to portrait orientation

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NoviceTester2Steps")

# ZipCode: keep the same value but store it as text (matches P2/Q2 text formatting)
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "363310"

# manufacturer / model / mobile_OS edits
$ws.Range("P2").NumberFormat = "@"
$ws.Range("P2").Value = "S"

$ws.Range("Q2").NumberFormat = "@"
$ws.Range("Q2").Value = "GG"

$ws.Range("R2").Value = "A"

# Move the sheet selection to S2
$ws.Activate()
$ws.Range("S2").Select()

# Switch the page to portrait orientation
$ws.PageSetup.Orientation = 1
